# Update countries & provincias Spain
# Applies the 11-Abril-2020 00:52 data refresh to the "Pais" sheet:
#  - Re-sorted country rows (Guinea/Isla de Man, Fiyi/Laos,
#    San Cristobal y Nieves/San Vicente y las Granadinas/Suazilandia/Seychelles)
#    now carry different country names at the same row positions.
#  - Updated case counters for the affected rows.
#  - Updated "last refreshed" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp header ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 00:52"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 501272
$ws.Range("C4").Value = 32706
$ws.Range("E4").Value = 455369
$ws.Range("G4").Value = 1973
$ws.Range("H4").Value = 18664

# --- Chequia (row 33) ---
$ws.Range("B33").Value = 5732
$ws.Range("C33").Value = 163
$ws.Range("E33").Value = 5267

# --- Argentina (row 56) ---
$ws.Range("B56").Value = 1975
$ws.Range("C56").Value = 81
$ws.Range("E56").Value = 1518
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 82

# --- Row 114: was "Isla de Man", now "Guinea" ---
$ws.Range("A114").Value = "Guinea"
$ws.Range("B114").Value = 212
$ws.Range("C114").Value = 18
$ws.Range("D114").Value = 15
$ws.Range("E114").Value = 197
$ws.Range("F114").Value = 0
$ws.Range("H114").Value = 0

# --- Row 115: was "Guinea", now "Isla de Man" ---
$ws.Range("A115").Value = "Isla de Man"
$ws.Range("B115").Value = 201
$ws.Range("C115").Value = 11
$ws.Range("H115").Value = 1

# --- Bahamas (row 150) ---
$ws.Range("B150").Value = 42
$ws.Range("C150").Value = 1
$ws.Range("E150").Value = 29

# --- Row 174: was "Laos", now "Fiyi" ---
$ws.Range("A174").Value = "Fiyi"
$ws.Range("C174").Value = 1

# --- Row 175: was "Fiyi", now "Laos" ---
$ws.Range("A175").Value = "Laos"
$ws.Range("C175").Value = 0

# --- Row 184: was "San Vicente y las Granadinas", now "San Cristobal y Nieves" ---
$ws.Range("A184").Value = "San Cristobal y Nieves"
$ws.Range("C184").Value = 1
$ws.Range("D184").Value = 0
$ws.Range("E184").Value = 12

# --- Row 185: was "Suazilandia", now "San Vicente y las Granadinas" ---
$ws.Range("A185").Value = "San Vicente y las Granadinas"
$ws.Range("D185").Value = 1
$ws.Range("E185").Value = 11

# --- Row 186: was "Seychelles", now "Suazilandia" ---
$ws.Range("A186").Value = "Suazilandia"
$ws.Range("B186").Value = 12
$ws.Range("D186").Value = 7
$ws.Range("E186").Value = 5

# --- Row 187: was "San Cristobal y Nieves", now "Seychelles" ---
$ws.Range("A187").Value = "Seychelles"
